# Applies the "updated parts list, including pricing" edit:
#  - case_1 (sheet1) row 10 (KY-016 indicator LED): supplier changes from
#    Amazon to a new "Banggood" supplier, price drops from 17.99 to 4.27,
#    and the notes cell gets a new note about shipping delay from China.
#  - The SUM total in C13 recalculates automatically from 86.11 to 72.39.
#  - Selection / active-tab bookkeeping: case_1 becomes the active sheet &
#    tab, case_10 is no longer the selected tab.

$wb = $excel.ActiveWorkbook

$wsCase1  = $wb.Worksheets.Item("case_1")
$wsCase10 = $wb.Worksheets.Item("case_10")

# --- Update parts list content on case_1 ---

# B10: Supplier -> Banggood (new shared string)
$wsCase1.Range("B10").Value = "Banggood"

# C10: Cost each -> 4.27
$wsCase1.Range("C10").Value = 4.27

# D10: Notes & alternatives -> new note about China shipping delay
$wsCase1.Range("D10").Value = "ALLOW 3+ WEEKS TO ARRIVE FROM CHINA. Due to COVID, it's hard to find these in the US. You can get it quickly from Amazon B07KJYR8K1, but costs `$18."

# C13 holds =SUM(C2:C12) already and will recalc automatically.

# --- Update sheet selections / active tab ---
# Order matters: selecting a range on a sheet also activates that sheet,
# so set case_10's selection first, then finish on case_1 so it ends up
# as the active/selected tab.

# case_10 selection moves to A12 (no longer the active tab afterwards).
$wsCase10.Range("A12").Select()

# case_1 becomes the selected/active tab with D6 as the active cell.
$wsCase1.Range("D6").Select()
